# Sync file from Google Drive
# Refresh the bus-arrival snapshot pulled into each NextBus* sheet:
# updated ETA timestamps (col F), recomputed MinutesToArrival (col O),
# and on NextBus2/NextBus3 a new bus (NextBus 151 -> Hougang Ctrl Int,
# service 74 already tracked) pushes the sheet from 7 to 8 data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# NextBus1: same 8 bus rows as before, only the polled ETA / minutes-to-
# arrival values changed.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("NextBus1")

$ws1.Cells.Item(2, 6).Value = 45688.60606481481
$ws1.Cells.Item(2, 15).Value = 7

$ws1.Cells.Item(3, 6).Value = 45688.61067129629
$ws1.Cells.Item(3, 15).Value = 14

$ws1.Cells.Item(4, 6).Value = 45688.62363425926
$ws1.Cells.Item(4, 15).Value = 32

$ws1.Cells.Item(5, 6).Value = 45688.61738425926
$ws1.Cells.Item(5, 15).Value = 23

$ws1.Cells.Item(6, 6).Value = 45688.6203125
$ws1.Cells.Item(6, 15).Value = 28

$ws1.Cells.Item(7, 15).Value = 21

$ws1.Cells.Item(8, 6).Value = 45688.60832175926
$ws1.Cells.Item(8, 15).Value = 10

# ---------------------------------------------------------------------
# NextBus2 and NextBus3 share the same refreshed snapshot: rows 2-5 keep
# their bus identity (only ETA/minutes move), row 6 turns into the new
# "Bt Batok Int" / SMRT 61 entry, and the two Hougang Ctrl Int buses
# (151 then 74) shift down into rows 7 and 8 - so a row must be inserted
# before rewriting the bottom of the table.
# ---------------------------------------------------------------------
$finalRows = @(
    @(2, @("NextBus3", 52, 53009, "Bishan Int", "SBST", 45688.61621527778, 53009, "WAB", "SEA", 1, 28009, "DD", 12101, "Ngee Ann Poly", 22)),
    @(3, @("NextBus3", 184, 44989, "Gali Batu Ter", "SMRT", 45688.61872685186, 44989, "WAB", "SEA", 1, 44989, "DD", 12101, "Ngee Ann Poly", 25)),
    @(4, @("NextBus3", 75, 44989, "Gali Batu Ter", "SMRT", 45688.62981481481, 44989, "WAB", "SEA", 0, 10009, "SD", 12101, "Ngee Ann Poly", 41)),
    @(5, @("NextBus3", 154, 82009, "Eunos Int", "SBST", 45688.6275, 82009, "WAB", "SEA", 1, 22009, "DD", 12101, "Ngee Ann Poly", 38)),
    @(6, @("NextBus3", 61, 43009, "Bt Batok Int", "SMRT", 45688.63138888889, 43009, "WAB", "SEA", 1, 82009, "SD", 12101, "Ngee Ann Poly", 44)),
    @(7, @("NextBus3", 151, 64009, "Hougang Ctrl Int", "SBST", 45688.62586805555, 64009, "WAB", "SEA", 0, 16009, "DD", 12101, "Ngee Ann Poly", 36)),
    @(8, @("NextBus3", 74, 64009, "Hougang Ctrl Int", "SBST", 45688.61626157408, 64009, "WAB", "SEA", 1, 11379, "DD", 12101, "Ngee Ann Poly", 22))
)

foreach ($sheetName in @("NextBus2", "NextBus3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Make room for the extra row (old row 7 slides down to row 8).
    $ws.Rows.Item(7).Insert()

    foreach ($entry in $finalRows) {
        $r = $entry[0]
        $vals = $entry[1]
        for ($i = 0; $i -lt $vals.Count; $i++) {
            $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
        }
    }
}
